$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture special row heights keyed by column A value (before any changes)
$specialHeights = @{}
for ($r = 2; $r -le 209; $r++) {
    $h = $ws.Cells.Item($r, 1).EntireRow.RowHeight
    if ($h -ne 12.75) {
        $a = $ws.Cells.Item($r, 1).Value2
        $specialHeights[$a] = $h
    }
}

# Add new row
$newRow = $ws.UsedRange.Rows.Count + 1
$ws.Cells.Item($newRow, 1).Value = "http://purl.obolibrary.org/obo/OBI_0000834"
$ws.Cells.Item($newRow, 2).Value = "high molecular weight DNA extract"
$ws.Cells.Item($newRow, 3).Value = "y"
$ws.Cells.Item($newRow, 4).Value = "genomic_DNA"

$ws.Range("C209").Copy()
$ws.Range("C210").PasteSpecial(-4122)

$rng = $ws.Range("A2:E210")
$rng.Sort($ws.Range("A2:A210"))

# Reset all data rows to default height, then reapply special ones at new positions
for ($r = 2; $r -le 210; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($specialHeights.ContainsKey($a)) {
        $ws.Rows($r).RowHeight = $specialHeights[$a]
    } else {
        $ws.Rows($r).RowHeight = 12.75
    }
}

Write-Host "Row20 height after fix:" $ws.Rows(20).RowHeight
Write-Host "Row11 height after fix:" $ws.Rows(11).RowHeight
